$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.357.35"
$ws.Range("E2").Value = "  +3.69%  "
$ws.Range("D3").Value = "3.245.99"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.86%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("D9").Value = "3.245.30"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("E10").Value = "  +5.73%  "
$ws.Range("E11").Value = "  +2.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.415"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.09%  "
$ws.Range("D13").Value = "3.808.31"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("E15").Value = "  +4.72%  "
$ws.Range("D16").Value = "67.343.01"
$ws.Range("E16").Value = "  +3.80%  "
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "3.241.05"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("E19").Value = "  +1.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("E30").Value = "  +8.19%  "
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +5.23%  "
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.34%  "
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.846"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  +4.20%  "
$ws.Range("E40").Value = "  +13.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "358.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.30%  "
$ws.Range("D45").Value = "2.716.95"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0681"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.88%  "
$ws.Range("E51").Value = "  -1.04%  "
